$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.703.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.56"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -1.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4865"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3790"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07336"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9159"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.57"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.877.41"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.472"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.609"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.92"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008807"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.741.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.120"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.102.98"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.905"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.66"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.150"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.84"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.902"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08907"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7656"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.644"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.517"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.093"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05270"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5467"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.917"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.456"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.92"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.60"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4789"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06054"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.86%  "
